$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# A1 holds a date serial value; update it to the new date (serial 45436 = 2024-05-24)
$ws.Range("A1").Value = (Get-Date -Year 2024 -Month 5 -Day 24 -Hour 0 -Minute 0 -Second 0)

# Update prices in column D
$ws.Range("D33").Value = 457
$ws.Range("D34").Value = 480
$ws.Range("D35").Value = 562
